$d = $word.ActiveDocument

# 1. Tweak wording: "of the sbt tool and" -> "of sbt and"
$d.Content.Find.Execute("of the sbt tool and", $true, $false, $false, $false, $false, $true, 1, $false, "of sbt and", 2) | Out-Null

# 2. Remove the first of the run of empty paragraphs that follows the
#    "... for an introduction in sbt. " paragraph (right before the
#    Bibliography heading).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "for an introduction in sbt\.") {
        $target = $p.Next()
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}
